$d = $word.ActiveDocument

# Disable "smart quotes" / autoformat-as-you-type so literal straight
# double quotes in the replacement text are not converted to curly quotes.
$word.Options.AutoFormatReplaceQuotes = $false
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

function Set-RangeText($old, $new) {
    # Locate the text with a throw-away Range (Content.Duplicate) and then
    # assign .Text directly on the located Range. This performs a plain
    # text substitution without Word's "replace as you type" autocorrect
    # (e.g. smart-quote conversion) that Find.Execute's ReplaceWith applies.
    $r = $d.Content.Duplicate
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
        return
    }
    $r.Text = $new
}

# --- Unique whole-string replacements -------------------------------------
Set-RangeText "2445987" "2487097"
Set-RangeText "генерального директора Котлярчука О. Е." "капитана Бахтин Ю. Г."
Set-RangeText "  Устава" "  Кодекса торгового мореплавания (КТМ РФ) "
Set-RangeText 'Рассмотрение технической документации "Грузовая марка" № 5234-234234-23 на т/х "МУРМАН 1" РС 091052' 'Рассмотрение технической документации "Переоборудование МО № 234-546-ИИ" на тх "ВОЛГА" РС 940330 / Review of technical documentation "" on mv VOLGA'
Set-RangeText "Дворцовая набережная, 8, Санкт-Петербург 191186" "Дворцовая набережная, 8, Санкт-Петербург  191186"
Set-RangeText "Молочинского, д. 4, Калининград 236023" "Молочинского, д. 4, Калининград  236023"
Set-RangeText "О. Е. Котлярчук" "Ю. Г. Бахтин"

# "ул. Карла Маркса..." appears twice and both occurrences change
# identically per the diff; replace each occurrence in turn.
Set-RangeText "ул. Карла Маркса, д. 19, Мурманск, Мурманская область 193025" "ул. Карла Маркса, д. 19, Мурманск  193025"
Set-RangeText "ул. Карла Маркса, д. 19, Мурманск, Мурманская область 193025" "ул. Карла Маркса, д. 19, Мурманск  193025"

# --- The lone "01" day-of-month field --------------------------------------
# "01" is also a substring of many unrelated numbers elsewhere in the
# document (account numbers, phone numbers, another date "09.01.2024"), so
# a document-wide Find for "01" is unsafe. The target is a single table
# cell holding exactly "01" (the day in "«01» мая ..."); locate that exact
# cell through the Tables collection and rewrite only its contents.
$t = $d.Tables.Item(1)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    for ($c = 1; $c -le $row.Cells.Count; $c++) {
        $cell = $row.Cells.Item($c)
        $cellRange = $cell.Range
        $txt = $cellRange.Text
        $txt = $txt -replace "`r", "" -replace "`a", ""
        if ($txt -eq "01") {
            $cellRange.Text = "03"
        }
    }
}
